$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.08106466666666667
$ws.Range("H2").Value = 0.243194
$ws.Range("I2").Value = 0.01252465659474717
$ws.Range("J2").Value = 0.01252465659474717
$ws.Range("M2").Value = 10.37574666666667
$ws.Range("N2").Value = 31.12724
$ws.Range("O2").Value = 0.8643482197679554
$ws.Range("P2").Value = 0.8643482197679554
$ws.Range("Q2").Value = 0.8411064449511112
$ws.Range("R2").Value = 7.56995800456
$ws.Range("S2").Value = 0.0108256646308747
$ws.Range("T2").Value = 0.0108256646308747

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.08106466666666667
$ws.Range("H3").Value = 0.243194
$ws.Range("I3").Value = 0.01252465659474717
$ws.Range("J3").Value = 0.01252465659474717
$ws.Range("O3").Value = 0.09618113591146868
$ws.Range("P3").Value = 0.09618113591146869
$ws.Range("Q3").Value = 0.09359488623644445
$ws.Range("R3").Value = 0.8423539761280001
$ws.Range("S3").Value = 0.00120463569818385
$ws.Range("T3").Value = 0.00120463569818385

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.08106466666666667
$ws.Range("H4").Value = 0.243194
$ws.Range("I4").Value = 0.01252465659474717
$ws.Range("J4").Value = 0.01252465659474717
$ws.Range("M4").Value = 0.4738106666666667
$ws.Range("N4").Value = 1.421432
$ws.Range("O4").Value = 0.03947064432057595
$ws.Range("P4").Value = 0.03947064432057595
$ws.Range("Q4").Value = 0.03840930375644445
$ws.Range("R4").Value = 0.345683733808
$ws.Range("S4").Value = 0.0004943562656886214
$ws.Range("T4").Value = 0.0004943562656886214

# Row 5
$ws.Range("I5").Value = 0.926055528343168
$ws.Range("J5").Value = 0.926055528343168
$ws.Range("M5").Value = 10.37574666666667
$ws.Range("N5").Value = 31.12724
$ws.Range("O5").Value = 0.8643482197679554
$ws.Range("P5").Value = 0.8643482197679554
$ws.Range("Q5").Value = 62.19022991805777
$ws.Range("R5").Value = 559.7120692625199
$ws.Range("S5").Value = 0.8004344473296906
$ws.Range("T5").Value = 0.8004344473296906

# Row 6
$ws.Range("I6").Value = 0.926055528343168
$ws.Range("J6").Value = 0.926055528343168
$ws.Range("O6").Value = 0.09618113591146868
$ws.Range("P6").Value = 0.09618113591146869
$ws.Range("S6").Value = 0.08906907263314118
$ws.Range("T6").Value = 0.08906907263314119

# Row 7
$ws.Range("I7").Value = 0.926055528343168
$ws.Range("J7").Value = 0.926055528343168
$ws.Range("M7").Value = 0.4738106666666667
$ws.Range("N7").Value = 1.421432
$ws.Range("O7").Value = 0.03947064432057595
$ws.Range("P7").Value = 0.03947064432057595
$ws.Range("Q7").Value = 2.839930006415111
$ws.Range("R7").Value = 25.559370057736
$ws.Range("S7").Value = 0.03655200838033622
$ws.Range("T7").Value = 0.03655200838033622

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.397534
$ws.Range("H8").Value = 1.192602
$ws.Range("I8").Value = 0.06141981506208485
$ws.Range("J8").Value = 0.06141981506208484
$ws.Range("M8").Value = 10.37574666666667
$ws.Range("N8").Value = 31.12724
$ws.Range("O8").Value = 0.8643482197679554
$ws.Range("P8").Value = 0.8643482197679554
$ws.Range("Q8").Value = 4.124712075386666
$ws.Range("R8").Value = 37.12240867847999
$ws.Range("S8").Value = 0.05308810780739009
$ws.Range("T8").Value = 0.05308810780739008

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.397534
$ws.Range("H9").Value = 1.192602
$ws.Range("I9").Value = 0.06141981506208485
$ws.Range("J9").Value = 0.06141981506208484
$ws.Range("O9").Value = 0.09618113591146868
$ws.Range("P9").Value = 0.09618113591146869
$ws.Range("Q9").Value = 0.4589810954026666
$ws.Range("R9").Value = 4.130829858624
$ws.Range("S9").Value = 0.005907427580143653
$ws.Range("T9").Value = 0.005907427580143654

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.397534
$ws.Range("H10").Value = 1.192602
$ws.Range("I10").Value = 0.06141981506208485
$ws.Range("J10").Value = 0.06141981506208484
$ws.Range("M10").Value = 0.4738106666666667
$ws.Range("N10").Value = 1.421432
$ws.Range("O10").Value = 0.03947064432057595
$ws.Range("P10").Value = 0.03947064432057595
$ws.Range("Q10").Value = 0.1883558495626667
$ws.Range("R10").Value = 1.695202646064
$ws.Range("S10").Value = 0.002424279674551104
$ws.Range("T10").Value = 0.002424279674551104
